# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) and "最低票价" (G column) values
# across the "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 141
$ws1.Range("F4").Value = 1766
$ws1.Range("F6").Value = 1037
$ws1.Range("F7").Value = 2192
$ws1.Range("G7").Value = 70
$ws1.Range("F8").Value = 2103
$ws1.Range("F9").Value = 1103
$ws1.Range("F10").Value = 603
$ws1.Range("F12").Value = 1667
$ws1.Range("F13").Value = 394
$ws1.Range("F17").Value = 203
$ws1.Range("F18").Value = 1576
$ws1.Range("F19").Value = 628
$ws1.Range("F20").Value = 712
$ws1.Range("F21").Value = 601
$ws1.Range("F22").Value = 12216
$ws1.Range("F23").Value = 12252
$ws1.Range("F24").Value = 909
$ws1.Range("F25").Value = 698
$ws1.Range("F29").Value = 360
$ws1.Range("F30").Value = 1919
$ws1.Range("F32").Value = 573

# ---- Sheet "演出" ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 21

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 141
$ws4.Range("F5").Value = 1766
$ws4.Range("F7").Value = 1037
$ws4.Range("F8").Value = 2192
$ws4.Range("G8").Value = 70
$ws4.Range("F9").Value = 2103
$ws4.Range("F10").Value = 1103
$ws4.Range("F11").Value = 603
$ws4.Range("F13").Value = 1667
$ws4.Range("F14").Value = 394
$ws4.Range("F21").Value = 203
$ws4.Range("F22").Value = 1576
$ws4.Range("F23").Value = 628
$ws4.Range("F24").Value = 712
$ws4.Range("F25").Value = 601
$ws4.Range("F26").Value = 12216
$ws4.Range("F27").Value = 12252
$ws4.Range("F28").Value = 909
$ws4.Range("F29").Value = 698
$ws4.Range("F33").Value = 360
$ws4.Range("F34").Value = 1919
$ws4.Range("F38").Value = 573
$ws4.Range("F39").Value = 21
